# "Added HEP B/D to Codebook"
#
# This script:
#  1. Adds codebook entries (Positive/Negative/Interdeterminate) for the
#     Hepatitis B/D related lab columns (LBXHBC, LBXBHS, LBDHBG, LBDHD)
#     to the "Codebook" sheet.
#  2. Fixes up the corresponding "attributeOf" mapping codes on the
#     "Dictionary Mapping" sheet for those same four columns, replacing
#     the placeholder "ncit:*" codes with the correct "obo:NCIT_*" codes,
#     and clears out the stray extra mapping cells that existed there.
#  3. Adds a missing Unit value for LBDNISI on "Dictionary Mapping".
#  4. Leaves the "Codebook" sheet active/selected as the last-touched sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Codebook sheet - add rows for LBXHBC, LBXBHS, LBDHBG, LBDHD
# ---------------------------------------------------------------------
$codebook = $wb.Worksheets.Item("Codebook")

$entries = @(
    @("LBXHBC", 1, "Positive"),
    @("LBXHBC", 2, "Negative"),
    @("LBXHBC", 3, "Interdeterminate"),
    @("LBXBHS", 1, "Positive"),
    @("LBXBHS", 2, "Negative"),
    @("LBXBHS", 3, "Interdeterminate"),
    @("LBDHBG", 1, "Positive"),
    @("LBDHBG", 2, "Negative"),
    @("LBDHBG", 3, "Interdeterminate"),
    @("LBDHD",  1, "Positive"),
    @("LBDHD",  2, "Negative"),
    @("LBDHD",  3, "Interdeterminate")
)

$row = 6
foreach ($entry in $entries) {
    $codebook.Range("A$row").Value = $entry[0]
    $codebook.Range("B$row").Value = $entry[1]
    $codebook.Range("C$row").Value = $entry[2]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Dictionary Mapping sheet - fix attributeOf codes for the same
#    four Hep B/D columns, and clean up stray cells in column I.
# ---------------------------------------------------------------------
$mapping = $wb.Worksheets.Item("Dictionary Mapping")

$mapping.Range("B25").Value = "obo: NCIT_C75678"
$mapping.Range("I25").ClearContents()

$mapping.Range("B26").Value = "obo: NCIT_C628795"
$mapping.Range("I26").ClearContents()

$mapping.Range("B27").Value = "obo: NCIT_C75678"
$mapping.Range("I27").ClearContents()

$mapping.Range("B28").Value = "obo:NCIT_C96664"
$mapping.Range("I28").ClearContents()

# Add missing Unit for LBDNISI (row 19)
$mapping.Range("D19").Value = "nhanes:00154"

# ---------------------------------------------------------------------
# 3. Make sure every worksheet has an explicit page setup, and leave
#    the Codebook sheet as the active / selected sheet & cell, matching
#    the saved workbook view state.
# ---------------------------------------------------------------------
foreach ($name in @("InfoSheet", "Prefixes", "Mapping Process", "Timeline", "New Concepts", "Dictionary Mapping")) {
    $sh = $wb.Worksheets.Item($name)
    $sh.PageSetup.Orientation = 1
}

$mapping.Activate()
$mapping.Range("A28").Select()

$codebook.Activate()
$excel.ActiveWindow.Zoom = 117
$codebook.Range("F36").Select()
